# Apply crypto price/volume updates per the diff (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.670.02"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "3.094.17"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.50"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.81"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.091.58"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.35"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -4.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.89"
$ws.Range("E14").Value = "  -4.72%  "

$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "3.612.87"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "66.660.67"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.97"
$ws.Range("E18").Value = "  +3.45%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.96"
$ws.Range("E19").Value = "  -2.88%  "

$ws.Range("D20").Value = "3.096.22"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.74"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.78"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.688"
$ws.Range("E23").Value = "  -3.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.34"

$ws.Range("E25").Value = "  -5.18%  "

$ws.Range("E26").Value = "  -3.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  -3.65%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  -4.93%  "

$ws.Range("E31").Value = "  -3.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.89"
$ws.Range("E32").Value = "  -2.93%  "

$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("D34").Value = "0.0₃0932"
$ws.Range("E34").Value = "  -2.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.87"
$ws.Range("E36").Value = "  +4.82%  "

$ws.Range("E37").Value = "  -6.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.942"
$ws.Range("E38").Value = "  -3.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.07"
$ws.Range("E39").Value = "  -2.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.309"
$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("E42").Value = "  -5.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.25"
$ws.Range("E43").Value = "  -3.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").Value = "2.778.00"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0346"
$ws.Range("E46").Value = "  -2.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "368.32"
$ws.Range("E47").Value = "  -4.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.33"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.46"
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.18"
$ws.Range("E51").Value = "  -2.36%  "
